$wb = $excel.ActiveWorkbook

$newHandoffDatetime = "2016-08-15 10:37:42"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06b4c25bc6dcf116ba277fc784512a0655008b37/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ed0a4373a1beed8bc06a2d29b68b3e43358b401/e2e/b.md."

# Overview sheet: row 3 corresponds to b.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $newHandoffDatetime

# zh-cn sheet: row 3 corresponds to b.md
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-15 10:37:38"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# de-de sheet: row 3 corresponds to b.md
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $newHandoffDatetime
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 40
